$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '43.726.66'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.40%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.312.15'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +4.02%  '

$ws.Range('E4').Value = '  +0.16%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '96.96'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +4.63%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '269.71'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.19%  '

$ws.Range('E7').Value = '  +0.17%  '

$ws.Range('E8').Value = '  +0.08%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.622'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.54%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '45.35'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.85%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0942'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.73%  '

$ws.Range('E12').Value = '  -3.14%  '

$ws.Range('E13').Value = '  +0.38%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.650.75'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.71%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.42'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.65%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.864'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +7.75%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.308.67'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +4.05%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '43.700.22'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.40%  '

$ws.Range('E19').Value = '  +5.06%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.36'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +5.53%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '72.55'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.02%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '238.96'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.80%  '

$ws.Range('E23').Value = '  -3.67%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.39'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.64%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.52'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.54%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.31'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.31%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '3.48'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.68%  '

$ws.Range('E29').Value = '  +0.41%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '22.36'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +7.20%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '37.89'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -9.31%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '174.29'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.97%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0897'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.56%  '

$ws.Range('E34').Value = '  -0.53%  '

$ws.Range('E35').Value = '  +2.78%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0361'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +3.03%  '

$ws.Range('E37').Value = '  -3.20%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.36'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.45%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.36'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -6.02%  '

$ws.Range('E40').Value = '  +10.46%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.36'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +9.27%  '

$ws.Range('E42').Value = '  +18.52%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '12.10'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.79%  '

$ws.Range('E44').Value = '  +9.28%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '61.77'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.46%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '5.30'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.32%  '

$ws.Range('E47').Value = '  +4.20%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '100.07'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.22%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.20'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.39%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.189'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +16.69%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.537.64'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.95%  '
